$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
$asv = $win.ActiveSheetView
Write-Host $asv
$asv | Get-Member
